# Update "Pais" country/cases table + timestamp to the 16:36 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Julio de 2020 a las 16:36"

# --- Row 4: Estados Unidos (no re-sort, still #1) ---
$ws.Range("B4").Value = 3098959
$ws.Range("C4").Value = 1875
$ws.Range("D4").Value = 1355675
$ws.Range("E4").Value = 1609254
$ws.Range("G4").Value = 58
$ws.Range("H4").Value = 134030

# --- Row 6: India (no re-sort, still #3) ---
$ws.Range("B6").Value = 753354
$ws.Range("C6").Value = 9873
$ws.Range("D6").Value = 464171
$ws.Range("E6").Value = 268380
$ws.Range("G6").Value = 150
$ws.Range("H6").Value = 20803

# --- Row 19: Alemania (no re-sort) ---
$ws.Range("B19").Value = 198512
$ws.Range("C19").Value = 157
$ws.Range("E19").Value = 6707
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 9105

# --- Rows 39-40: Filipinas overtakes Oman ---
# Row 39 becomes Filipinas with its updated totals.
$ws.Range("A39").Value = "Filipinas"
$ws.Range("B39").Value = 50359
$ws.Range("C39").Value = 2486
$ws.Range("D39").Value = 12588
$ws.Range("E39").Value = 36457
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 5
$ws.Range("H39").Value = 1314
# Row 40 becomes Oman, unchanged totals, now ranked below Filipinas.
$ws.Range("A40").Value = "Oman"
$ws.Range("B40").Value = 50207
$ws.Range("C40").Value = 1210
$ws.Range("D40").Value = 32005
$ws.Range("E40").Value = 17969
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 9
$ws.Range("H40").Value = 233

# --- Rows 62-63: Serbia overtakes Argelia ---
# Row 62 becomes Serbia with its updated totals.
$ws.Range("A62").Value = "Serbia"
$ws.Range("B62").Value = 17076
$ws.Range("C62").Value = 357
$ws.Range("D62").Value = 13366
$ws.Range("E62").Value = 3369
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 11
$ws.Range("H62").Value = 341
# Row 63 becomes Argelia, unchanged totals, now ranked below Serbia.
$ws.Range("A63").Value = "Argelia"
$ws.Range("B63").Value = 16879
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 12094
$ws.Range("E63").Value = 3817
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 968

# --- Row 71: Uzbekistan (no re-sort) ---
$ws.Range("B71").Value = 10982
$ws.Range("C71").Value = 312
$ws.Range("D71").Value = 6888
$ws.Range("E71").Value = 4052

# --- Row 85: Tayikistan (no re-sort) ---
$ws.Range("B85").Value = 6364
$ws.Range("C85").Value = 49
$ws.Range("D85").Value = 5011
$ws.Range("E85").Value = 1299
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 54

# --- Row 161: Vietnam (active/recovered reclassification only) ---
$ws.Range("D161").Value = 347
$ws.Range("E161").Value = 22

# --- Rows 209-210: Groenlandia / Islas Malvinas tie swap (values identical) ---
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"
